# Update export excel format
# 1. Insert a new first sheet "NOTICE et COMMENTAIRES" with a single informational cell.
# 2. Update the MAGISTRATS / FONCTIONNAIRES header cell text + formatting (wider column,
#    taller header row, wrap text) and refresh the remembered selection.

$wb = $excel.ActiveWorkbook

$magistrats = $wb.Worksheets.Item("MAGISTRATS")

# --- New "NOTICE et COMMENTAIRES" sheet, inserted before MAGISTRATS -----------------
$notice = $wb.Worksheets.Add($magistrats)
$notice.Name = "NOTICE et COMMENTAIRES"
$notice.Range("A1").Value = "#! FINISH"

# Re-resolve the other sheets by name now that the sheet collection has shifted.
$magistrats     = $wb.Worksheets.Item("MAGISTRATS")
$fonctionnaires = $wb.Worksheets.Item("FONCTIONNAIRES")

# --- MAGISTRATS sheet -----------------------------------------------------------------
$magistrats.Range("B1").Value = "#`` Export d'un référentiel de temps moyens A-JUST : `${name} `n(ce fichier peut être importé directement dans A-JUST)"
$magistrats.Range("B1:C1").WrapText = $true
$magistrats.Columns.Item(2).ColumnWidth = 50.83203125
$magistrats.Rows.Item(1).RowHeight = 59
$magistrats.Range("B58").Select()

# --- FONCTIONNAIRES sheet ---------------------------------------------------------------
$fonctionnaires.Range("B1").Value = "#`` Export d'un référentiel de temps moyens A-JUST : `${nameFonc} `n(ce fichier peut être importé directement dans A-JUST)"
$fonctionnaires.Range("B1:C1").WrapText = $true
$fonctionnaires.Columns.Item(2).ColumnWidth = 50.83203125
$fonctionnaires.Rows.Item(1).RowHeight = 60
$fonctionnaires.Range("B40").Select()

# Leave the new notice sheet active/selected, matching the authored workbook.
$notice.Activate()
